# teamdata.xlsx edit: convert season totals to per-game averages and
# rename the shooting-percentage header columns accordingly.
#
# Column layout (unchanged positions) is:
#   A Rk   B Team  C G   D MP   E FG   F FGA
#   G FG%->FGP   H 3P->ThreeP   I 3PA->ThreePA   J 3P%->ThreePP
#   K 2P->TwoP   L 2PA->TwoPA   M 2P%->TwoPP
#   N FT (unchanged)  O FTA (unchanged)  P FT%->FTP
#   Q ORB  R DRB  S TRB  T AST  U STL  V BLK  W TOV  X PF  Y PTS
#   Z Playoff_Wins  AA Playoff_Birth
#
# Every counting-stat column (everything except the percentage columns
# G/J/M/P, the Rk/Team/G columns, and the Z/AA playoff columns) is
# divided by column C (games played) to turn season totals into
# per-game rates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells (row 1) ---------------------------------
$ws.Range("G1").Value = "FGP"
$ws.Range("H1").Value = "ThreeP"
$ws.Range("I1").Value = "ThreePA"
$ws.Range("J1").Value = "ThreePP"
$ws.Range("K1").Value = "TwoP"
$ws.Range("L1").Value = "TwoPA"
$ws.Range("M1").Value = "TwoPP"
$ws.Range("P1").Value = "FTP"

# --- 2. Convert season totals to per-game averages ------------------
# Columns to divide by games-played (column C).
$statCols = @(4,5,6,8,9,11,12,14,15,17,18,19,20,21,22,23,24,25)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 31 }

for ($r = 2; $r -le $lastRow; $r++) {
    $games = $ws.Cells.Item($r, 3).Value2
    foreach ($col in $statCols) {
        $cell = $ws.Cells.Item($r, $col)
        $total = $cell.Value2
        if ($games -ne 0) {
            $cell.Value = $total / $games
        }
    }
}

# --- 3. Ensure the Playoff_Birth (AA) column has no blank cells -----
# (Dallas Mavericks / last row had an implicit blank meaning "0".)
for ($r = 2; $r -le $lastRow; $r++) {
    $aaCell = $ws.Cells.Item($r, 27)
    if ($aaCell.Value2 -eq $null) {
        $aaCell.Value = 0
    }
}

# --- 4. Restore the original selection / active cell ----------------
$ws.Range("P1").Select()
